# Add a new "Addition features" section with one new feature row
# ("Motion blur acceleration structure") to the workload sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")

# --- Row 28: section header "Addition features" ---------------------------
$b28 = $ws.Range("B28")
$b28.Value = "Addition features"
$b28.Borders.LineStyle = 1
$b28.Font.Bold = $true
$b28.HorizontalAlignment = -4108   # xlCenter

$c28 = $ws.Range("C28")
$c28.Borders.LineStyle = 1
$c28.NumberFormat = "0.00"
$c28.HorizontalAlignment = -4108   # xlCenter

$ws.Range("D28:G28").Borders.LineStyle = 1

# --- Row 29: new feature "Motion blur acceleration structure" -------------
$b29 = $ws.Range("B29")
$b29.Value = "Motion blur acceleration structure"
$b29.Borders.LineStyle = 1
$b29.HorizontalAlignment = -4108   # xlCenter

$c29 = $ws.Range("C29")
$c29.Borders.LineStyle = 1
$c29.NumberFormat = "0.00"
$c29.HorizontalAlignment = -4108   # xlCenter

$ws.Range("D29:G29").Borders.LineStyle = 1
$ws.Range("E29").Value = 100

# --- View bookkeeping: scroll down and move the selection ------------------
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E30").Select() | Out-Null
